$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.464.68'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '3.438.78'
$ws.Range("E3").Value = '  +1.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.49%  '
$ws.Range("D7").Value = '3.439.64'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.475'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.93%  '
$ws.Range("E10").Value = '  +0.69%  '
$ws.Range("E11").Value = '  +3.51%  '
$ws.Range("E12").Value = '  +2.05%  '
$ws.Range("D13").Value = '4.027.16'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.36%  '
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000173'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("D17").Value = '3.445.49'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").Value = '61.578.42'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.27'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.50%  '
$ws.Range("E20").Value = '  +3.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '394.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.54%  '
$ws.Range("E23").Value = '  +3.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.10'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.43%  '
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '3.587.72'
$ws.Range("E28").Value = '  +2.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.177'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.08%  '
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("E32").Value = '  -7.94%  '
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("E34").Value = '  +2.52%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  +3.24%  '
$ws.Range("D37").Value = '3.470.06'
$ws.Range("E37").Value = '  +2.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("E40").Value = '  +1.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '167.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0779'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.91%  '
$ws.Range("E44").Value = '  +4.06%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  +1.44%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.04%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("D49").Value = '2.600.27'
$ws.Range("E49").Value = '  +2.49%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("E51").Value = '  +2.59%  '
